# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 7316
    4  = 5636
    6  = 175
    10 = 88
    11 = 114
    12 = 206
    13 = 55
    15 = 362
    17 = 13
    18 = 1
    20 = 50
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
